$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "dfdsfsdf"
$ws.Range("B8").Value = 3
$ws.Range("C8").Value = 6
